# Adds a "loose" clean mode: cell A2 ("Robert Eshleman") is rewritten with
# leading/trailing whitespace plus an internal newline so the cleaning
# logic has real padding/line-breaks to exercise, and the cell is given
# wrap-text formatting (with a taller row) so the multi-line value renders.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New padded, multi-line value for the name cell.
$ws.Range("A2").Value = "    Robert`nEshleman    "

# Wrap the text so the embedded newline actually breaks the line, and
# bump the row height to fit both wrapped lines.
$ws.Range("A2").WrapText = $true
$ws.Rows.Item(2).RowHeight = 34
